$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells with new values (rows 86, 87, 90, 91, 92) ---
$ws.Range("F86").Value2 = 9.407253593377034
$ws.Range("Q86").Value2 = 31.83817342823687
$ws.Range("F87").Value2 = 5.822614627431119
$ws.Range("Q87").Value2 = 25.49896763766734
$ws.Range("C90").Value2 = 6.844717012085865
$ws.Range("E90").Value2 = 9.538919156451357
$ws.Range("F90").Value2 = 9.792327288220028
$ws.Range("K90").Value2 = 5.765605168853732
$ws.Range("L90").Value2 = 2.610188788384008
$ws.Range("M90").Value2 = 7.8313747381388
$ws.Range("P90").Value2 = 49.98061607552975
$ws.Range("Q90").Value2 = 32.40251607660404
$ws.Range("F91").Value2 = 7.833333333333333
$ws.Range("J91").Value2 = 9.887161667102284
$ws.Range("Q91").Value2 = 37.72049500043562
$ws.Range("F92").Value2 = 2.333333333333332
$ws.Range("Q92").Value2 = 32.33333333333333

# --- Append new rows 94-101 ---
$ws.Range("A94").NumberFormat = "@"
$ws.Range("A94").Value2 = "2025-02-24"
$ws.Range("A94").ClearFormats()
$ws.Range("B94").Value2 = "abs_activity"
$arr94 = New-Object 'object[,]' 1,15
$arr94[0,0] = 8.35475344542496
$arr94[0,1] = 0
$arr94[0,2] = 7.878390754910651
$arr94[0,3] = 9.310956088448684
$arr94[0,4] = 10
$arr94[0,5] = 10
$arr94[0,6] = 6.941346394792772
$arr94[0,7] = 10
$arr94[0,8] = -8.964677414461047
$arr94[0,9] = 6.763844676910046
$arr94[0,10] = 9.00656105060188
$arr94[0,11] = 0
$arr94[0,12] = 0
$arr94[0,13] = 33.21637423126921
$arr94[0,14] = 36.07480076535873
$ws.Range("C94:Q94").Value2 = $arr94

$ws.Range("A95").NumberFormat = "@"
$ws.Range("A95").Value2 = "2025-02-24"
$ws.Range("A95").ClearFormats()
$ws.Range("B95").Value2 = "rel_activity"
$arr95 = New-Object 'object[,]' 1,15
$arr95[0,0] = 5.971961929959944
$arr95[0,1] = 5
$arr95[0,2] = 0
$arr95[0,3] = 6.583333333333333
$arr95[0,4] = 7.751977984176126
$arr95[0,5] = 0
$arr95[0,6] = 5.062003968253968
$arr95[0,7] = 10
$arr95[0,8] = 0
$arr95[0,9] = 0
$arr95[0,10] = 5.434027777777778
$arr95[0,11] = 5
$arr95[0,12] = 5
$arr95[0,13] = 29.21997166016781
$arr95[0,14] = 26.58333333333333
$ws.Range("C95:Q95").Value2 = $arr95

$ws.Range("A96").NumberFormat = "@"
$ws.Range("A96").Value2 = "2025-02-24"
$ws.Range("A96").ClearFormats()
$ws.Range("B96").Value2 = "abs_sleep"
$arr96 = New-Object 'object[,]' 1,15
$arr96[0,0] = 10
$arr96[0,1] = 0
$arr96[0,2] = 9.666666666666666
$arr96[0,3] = 9.266666666666667
$arr96[0,4] = 10
$arr96[0,5] = 5.733333333333334
$arr96[0,6] = 10
$arr96[0,7] = 10
$arr96[0,8] = 9.566666666666666
$arr96[0,9] = 10
$arr96[0,10] = 10
$arr96[0,11] = 0
$arr96[0,12] = 0
$arr96[0,13] = 59.23333333333333
$arr96[0,14] = 35
$ws.Range("C96:Q96").Value2 = $arr96

$ws.Range("A97").NumberFormat = "@"
$ws.Range("A97").Value2 = "2025-02-24"
$ws.Range("A97").ClearFormats()
$ws.Range("B97").Value2 = "rel_sleep"
$arr97 = New-Object 'object[,]' 1,15
$arr97[0,0] = 9.852434077079108
$arr97[0,1] = 0
$arr97[0,2] = 0
$arr97[0,3] = 8.384657063256086
$arr97[0,4] = 8.489671610169495
$arr97[0,5] = 7.108288914824752
$arr97[0,6] = 9.99232158988257
$arr97[0,7] = 10
$arr97[0,8] = 0
$arr97[0,9] = 8.310295427942485
$arr97[0,10] = 0
$arr97[0,11] = 0
$arr97[0,12] = 0
$arr97[0,13] = 28.33442727713118
$arr97[0,14] = 33.80324140602332
$ws.Range("C97:Q97").Value2 = $arr97

$ws.Range("A98").NumberFormat = "@"
$ws.Range("A98").Value2 = "2025-02-25"
$ws.Range("A98").ClearFormats()
$ws.Range("B98").Value2 = "abs_activity"
$arr98 = New-Object 'object[,]' 1,15
$arr98[0,0] = 2.802331824079988
$arr98[0,1] = 0
$arr98[0,2] = 9.455073435052073
$arr98[0,3] = 0
$arr98[0,4] = 8.677223592853386
$arr98[0,5] = 9.505661273856012
$arr98[0,6] = 0
$arr98[0,7] = 9.431214239560413
$arr98[0,8] = 0
$arr98[0,9] = 8.744036274687891
$arr98[0,10] = 6.042290974649834
$arr98[0,11] = 0
$arr98[0,12] = 0
$arr98[0,13] = 26.97691982663528
$arr98[0,14] = 27.68091178810432
$ws.Range("C98:Q98").Value2 = $arr98

$ws.Range("A99").NumberFormat = "@"
$ws.Range("A99").Value2 = "2025-02-25"
$ws.Range("A99").ClearFormats()
$ws.Range("B99").Value2 = "rel_activity"
$arr99 = New-Object 'object[,]' 1,15
$arr99[0,0] = 0
$arr99[0,1] = 5
$arr99[0,2] = 6.456700823789431
$arr99[0,3] = 5.602271423654926
$arr99[0,4] = 6.891984864121087
$arr99[0,5] = 0
$arr99[0,6] = 0
$arr99[0,7] = 6.083333333333333
$arr99[0,8] = 0
$arr99[0,9] = 0
$arr99[0,10] = 0
$arr99[0,11] = 5
$arr99[0,12] = 5
$arr99[0,13] = 18.34868568791052
$arr99[0,14] = 21.68560475698826
$ws.Range("C99:Q99").Value2 = $arr99

$ws.Range("A100").NumberFormat = "@"
$ws.Range("A100").Value2 = "2025-02-25"
$ws.Range("A100").ClearFormats()
$ws.Range("B100").Value2 = "abs_sleep"
$arr100 = New-Object 'object[,]' 1,15
$arr100[0,0] = 8.533333333333333
$arr100[0,1] = 0
$arr100[0,2] = 10
$arr100[0,3] = 0
$arr100[0,4] = 8.566666666666666
$arr100[0,5] = 8.1
$arr100[0,6] = 0
$arr100[0,7] = 10
$arr100[0,8] = 0
$arr100[0,9] = 4.600000000000001
$arr100[0,10] = 9.166666666666666
$arr100[0,11] = 0
$arr100[0,12] = 0
$arr100[0,13] = 36.26666666666667
$arr100[0,14] = 22.7
$ws.Range("C100:Q100").Value2 = $arr100

$ws.Range("A101").NumberFormat = "@"
$ws.Range("A101").Value2 = "2025-02-25"
$ws.Range("A101").ClearFormats()
$ws.Range("B101").Value2 = "rel_sleep"
$arr101 = New-Object 'object[,]' 1,15
$arr101[0,0] = 0
$arr101[0,1] = 0
$arr101[0,2] = 7.119049981150419
$arr101[0,3] = 0
$arr101[0,4] = 0
$arr101[0,5] = 8.443852197663364
$arr101[0,6] = 0
$arr101[0,7] = 10
$arr101[0,8] = 0
$arr101[0,9] = 0
$arr101[0,10] = 0
$arr101[0,11] = 0
$arr101[0,12] = 0
$arr101[0,13] = 7.119049981150419
$arr101[0,14] = 18.44385219766336
$ws.Range("C101:Q101").Value2 = $arr101
